# Updated ITA model - 2025-08-04 23:53
#
# VEDA_Sets-Proc: extend a couple of process-name matching patterns with
# extra wildcard aliases, and fill in the And/Or (T_Pos_AndOr / T_Neg_AndOr)
# helper columns plus a SetDesc value that were left blank for a few rows.
# Also add a new "-*SMR" exclusion pattern for the nuclear set (so SMR
# process names are excluded from the generic nuclear bucket).

$wb = $excel.ActiveWorkbook

# Workbook no longer auto-updates external links on open.
try { $wb.UpdateLinks = 0 } catch { }

$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

# PSET_PN (col B): append new wildcard aliases to the existing patterns
$ws.Range("B3").Value = "ep_gas_combined_cycle*,ep_oil_combined_cycle*,CCGT*,*GasCC*"
$ws.Range("B7").Value = "ep_gas_gas_turbine*,ep_oil_gas_turbine*,gas turbine*,EN*CT*"

# Row 3 (CCGT set): fill in SetDesc (G) and the And/Or columns (H/I)
$ws.Range("G3").Value = "CCGT"
$ws.Range("H3").Value = "And"
$ws.Range("I3").Value = "Or"

# Row 7 (OCGT / Peaker set): fill in the And/Or columns (H/I)
$ws.Range("H7").Value = "And"
$ws.Range("I7").Value = "Or"

# Row 17 (Nuclear set): add a PSET_PN exclusion for SMR units, plus And/Or
$ws.Range("B17").Value = "-*SMR"
$ws.Range("H17").Value = "And"
$ws.Range("I17").Value = "Or"

$wb.Save()
